# Rename the header/footer logo images (inline pictures) so that the
# wp:docPr / pic:cNvPr "name" attribute matches the new filenames.
#
# Mapping of logical Word object -> physical header/footer part
# (discovered by probing the default/first-page link order in this
# single-section document):
#   Sections(1).Headers(1)  (wdHeaderFooterPrimary)   -> BTec logo,   id=3
#   Sections(1).Headers(2)  (wdHeaderFooterFirstPage)  -> BTec logo,   id=1
#   Sections(1).Footers(1)  (wdHeaderFooterPrimary)   -> Pearson logo, id=4
#   Sections(1).Footers(2)  (wdHeaderFooterFirstPage)  -> Pearson logo, id=2
#
# Diff summary: the BTec logo picture is renamed image1.jpg -> image2.jpg
# (in both headers), and the Pearson logo picture is renamed
# image2.png -> image1.png (in both footers).

$d = $word.ActiveDocument
$sec = $d.Sections.First

# --- Headers: BTec_Logo-Orange, image1.jpg -> image2.jpg ---
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists) {
        $shapes = $hdr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

# --- Footers: PearsonLogo.png, image2.png -> image1.png ---
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists) {
        $shapes = $ftr.Range.InlineShapes
        for ($j = 1; $j -le $shapes.Count; $j++) {
            $shp = $shapes.Item($j)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
